$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices corresponding to B,C,D,E,F,I,J,K,L,M,N (G and H are unchanged/absent)
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

# New vm_pu values for rows 2..25 (one array per row, aligned with $cols)
$data = @(
    @(1.02, 1.029392447873921, 1.034854849791811, 1.033027458829641, 1.040262791202527, 1.03741252832598, 1.034539476242356, 1.037653127440003, 1.03583099068737, 1.043045637078266, 1.015505892564679),
    @(1.02, 1.030266124941358, 1.035530257684425, 1.033848708211137, 1.041402718879407, 1.037668335386179, 1.035054553785757, 1.038138165883869, 1.036461101932312, 1.043995085100992, 1.015678289203502),
    @(1.02, 1.030831781883657, 1.035967481272318, 1.034380784426401, 1.042140988850831, 1.037832677902401, 1.035387536900581, 1.038451528209544, 1.036868850863004, 1.044609522175739, 1.015789691298532),
    @(1.02, 1.031069661906005, 1.036151334133309, 1.034604629058861, 1.042451515378038, 1.03790148404017, 1.035527448784586, 1.038583147531222, 1.037040273426908, 1.044867850452436, 1.015836488480076),
    @(1.02, 1.031109607543419, 1.03618220638362, 1.034642222889666, 1.042503663394961, 1.037913020255767, 1.03555093622896, 1.038605240037086, 1.037069056287909, 1.04491122600557, 1.015844343796775),
    @(1.02, 1.030834960143567, 1.035969937751743, 1.034383774823112, 1.042145137501684, 1.037833598407843, 1.0353894067027, 1.038453287378933, 1.036871141401675, 1.044612973896915, 1.015790316747656),
    @(1.02, 1.029687642293268, 1.035083067025428, 1.033304863774452, 1.040647898144869, 1.037499223884528, 1.034713612066261, 1.037817149367203, 1.036043933974875, 1.043366490532173, 1.015564185685095),
    @(1.02, 1.027668490293716, 1.033521796540465, 1.031408896592052, 1.03801464321275, 1.036900986257551, 1.033520465225905, 1.036692479786062, 1.03458652423704, 1.041170662940938, 1.015164579352531),
    @(1.02, 1.026324175939812, 1.032482041093079, 1.030148499994683, 1.03626256532781, 1.03649612353197, 1.032723529876015, 1.035940256188257, 1.033615136396569, 1.039707226049465, 1.01489743201704),
    @(1.02, 1.025742509078549, 1.032032090420516, 1.029603599893673, 1.035504711086565, 1.036319387919375, 1.032378101447146, 1.03561396627639, 1.033194579075427, 1.039073651732014, 1.014781582085267),
    @(1.02, 1.02552651766613, 1.031865000582346, 1.029401330198735, 1.035223331616546, 1.036253526477726, 1.03224974212184, 1.03549268255977, 1.033038375352082, 1.038838329780216, 1.014738524531836),
    @(1.02, 1.025572845575557, 1.031900840004559, 1.029444711810793, 1.035283682975707, 1.03626766365519, 1.032277277965557, 1.03551870216684, 1.033071881113247, 1.038888806421151, 1.014747761682489),
    @(1.02, 1.025724653818214, 1.032018277860801, 1.029586877548831, 1.035481449695986, 1.03631394815411, 1.032367492278, 1.035603942658004, 1.033181667017485, 1.039054199623357, 1.014778023457533),
    @(1.02, 1.025818196602049, 1.03209064078675, 1.029674487877781, 1.035603316380802, 1.036342437223351, 1.032423069449197, 1.035656450884446, 1.033249311089496, 1.03915610592851, 1.014796665335824),
    @(1.02, 1.026362788395891, 1.032511908640359, 1.030184681446026, 1.036312878628859, 1.036507822860333, 1.032746447512355, 1.035961898992745, 1.033643048748603, 1.039749276488826, 1.014905116969749),
    @(1.02, 1.026704512257717, 1.032776232096321, 1.030504943518483, 1.036758184410859, 1.036611183090306, 1.032949200830767, 1.036153345939403, 1.033890046736904, 1.040121384568808, 1.014973099621372),
    @(1.02, 1.026903875442669, 1.032930433491169, 1.03069182985402, 1.037018001601353, 1.036671333609791, 1.033067429572975, 1.036264958445734, 1.034034122210801, 1.04033843888966, 1.015012736002724),
    @(1.02, 1.026971860136836, 1.032983016552524, 1.030755567285373, 1.03710660575312, 1.036691820002414, 1.033107736755462, 1.036303006027663, 1.034083249165811, 1.040412450489293, 1.01502624813388),
    @(1.02, 1.026667844228467, 1.032747869993075, 1.030470573846459, 1.036710399260422, 1.03660010776086, 1.032927450812904, 1.036132811203426, 1.033863545590983, 1.040081459860243, 1.014965807456384),
    @(1.02, 1.025679948257194, 1.031983694173064, 1.029545009670738, 1.035423208975782, 1.036300324425419, 1.032340927835908, 1.035578843800765, 1.033149337521411, 1.039005495010887, 1.014769112826784),
    @(1.02, 1.025059198618262, 1.031503469410103, 1.028963826190869, 1.034614603431989, 1.036110601042649, 1.031971858095529, 1.035230050473132, 1.03270034484805, 1.038329084492064, 1.014645294393932),
    @(1.02, 1.025388233435586, 1.031758022196168, 1.029271850518624, 1.035043193949986, 1.036211294156013, 1.03216753707279, 1.035414998773982, 1.032938358509582, 1.038687653782549, 1.014710946878443),
    @(1.02, 1.026684412804554, 1.032760685526509, 1.030486103765838, 1.036731991072313, 1.03660511265117, 1.032937278814372, 1.036142090137517, 1.033875520300489, 1.040099500086597, 1.014969102523635),
    @(1.02, 1.028190178916266, 1.033925235617418, 1.031898423944132, 1.038694800112738, 1.037056711517877, 1.033829190999873, 1.036983668432227, 1.034963265467311, 1.041738259566917, 1.015268019585527)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $excelRow = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($excelRow, $cols[$j]).Value = $rowValues[$j]
    }
}

Write-Output "Updated vm_pu values for rows 2-25 (B,C,D,E,F,I,J,K,L,M,N)."
